$wb = $excel.ActiveWorkbook

# Map of sheet name -> { row -> newValue } for column F ("想去人数")
$changes = @{
    "展览" = @{
        2 = 226; 3 = 754; 5 = 2254; 6 = 1341; 7 = 801; 8 = 110; 9 = 26; 10 = 2919;
        11 = 27; 16 = 110; 17 = 103; 18 = 980; 19 = 980; 20 = 115; 22 = 133; 24 = 168;
        25 = 624; 27 = 298; 28 = 30; 30 = 998; 31 = 4956; 32 = 434; 33 = 196; 34 = 92
    }
    "演出" = @{
        6 = 400; 11 = 189; 20 = 33; 21 = 296; 24 = 367; 26 = 624; 33 = 272; 38 = 737; 39 = 36
    }
    "本地生活" = @{
        5 = 398; 6 = 379
    }
    "全部类型" = @{
        3 = 226; 4 = 398; 6 = 754; 8 = 400; 10 = 2254; 11 = 1341; 12 = 801; 13 = 110;
        15 = 189; 16 = 26; 17 = 2919; 18 = 27; 24 = 379; 26 = 980; 27 = 980; 28 = 115;
        31 = 296; 32 = 133; 33 = 168; 36 = 624; 38 = 367; 39 = 624; 40 = 298; 42 = 30;
        44 = 998; 45 = 4956; 47 = 434; 48 = 196; 49 = 737
    }
}

foreach ($sheetName in $changes.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowChanges = $changes[$sheetName]
    foreach ($row in $rowChanges.Keys) {
        $ws.Cells.Item($row, 6).Value = $rowChanges[$row]
    }
}
